# Insert a new row above row 3 (shifts existing rows 3-7 down to rows 4-8)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("3:3").Insert()

# Populate the new row 3 with the new weekly record.
# Columns A, B, C, E, F, G, H, I, N, O, Q, R mirror the constant values
# used throughout this data set (and match what used to be in row 3
# before the insert shifted it down to row 4).
$ws.Range("A3").Value = 11
$ws.Range("B3").Value = "Vega Monumental Concepción"
$ws.Range("C3").Value = "Bíobío"
$ws.Range("D3").Value = 44519
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 300000000
$ws.Range("G3").Value = "Espárragos"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 250
$ws.Range("K3").Value = 1200
$ws.Range("L3").Value = 1300
$ws.Range("M3").Value = 1240
$ws.Range("N3").Value = "$/kilo"
$ws.Range("O3").Value = "Provincia de Linares"
$ws.Range("P3").Value = 1240
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = "Hortaliza"
